# Generate Report for Handoff
# Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the b753c5eb-... row (row 7) across the Overview, zh-cn and
# de-de sheets, as produced by a fresh CI report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-22 18:44:39"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-22 18:44:34"

# de-de sheet: column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-22 18:44:39"
